$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: JST connector changed from 1.25mm (GH) pitch to 2.0mm (PH) pitch part.
$ws.Range("C31").Value = "S3B-PH-SM4-TB"
$ws.Range("B31").Value = "CONN HEADER PH SIDE 3POS 2MM SMD"
$ws.Range("D31").Value = "http://www.jst-mfg.com/product/pdf/eng/ePH.pdf"
$ws.Range("E31").Value = 0.96

# Update the selected cell to match the saved view state.
$ws.Range("E32").Select()
